$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before DL, shifting DL:MN -> DM:MO
$ws.Range("DL1").EntireColumn.Insert()

# Set header for the newly inserted column
$ws.Range("DL1").Value = "DemonstrationProjectIdentifier"

# Clear the OtherSubscriberRelationshipType value in row 2 (was "SELF")
$ws.Range("AT2").Clear()

# Update the Id column value for all data rows
$ws.Range("A2").Value = "690148897e79911955eafc9a"
$ws.Range("A3").Value = "690148897e79911955eafc9a"
$ws.Range("A4").Value = "690148897e79911955eafc9a"
